# Weekly refresh of the Cilantro / Lo Valledor sheet:
# two new daily records are prepended to the existing price history
# (rows 799-800), which pushes every subsequent row down by two and
# makes the last two rows of the former table (what used to be rows
# 836-837) become the new last two rows (838-839).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 799 - this shifts
# old rows 799..837 down to 801..839 and keeps every other row (and the
# header) untouched.
$ws.Rows.Item(799).Insert()
$ws.Rows.Item(799).Insert()

# --- New row 799 --------------------------------------------------------
$ws.Cells.Item(799, 1).Value  = 6
$ws.Cells.Item(799, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(799, 3).Value  = "Metropolitana"
$ws.Cells.Item(799, 4).Value  = 44706
$ws.Cells.Item(799, 5).Value  = 13
$ws.Cells.Item(799, 6).Value  = 100112040
$ws.Cells.Item(799, 7).Value  = "Cilantro"
$ws.Cells.Item(799, 8).Value  = "Sin especificar"
$ws.Cells.Item(799, 9).Value  = "Primera"
$ws.Cells.Item(799, 10).Value = 570
$ws.Cells.Item(799, 11).Value = 4500
$ws.Cells.Item(799, 12).Value = 5000
$ws.Cells.Item(799, 13).Value = 4719
$ws.Cells.Item(799, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(799, 15).Value = "Región Metropolitana"
$ws.Cells.Item(799, 16).Value = 131
$ws.Cells.Item(799, 17).Value = 36
$ws.Cells.Item(799, 18).Value = "Hortaliza"

# --- New row 800 --------------------------------------------------------
$ws.Cells.Item(800, 1).Value  = 6
$ws.Cells.Item(800, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(800, 3).Value  = "Metropolitana"
$ws.Cells.Item(800, 4).Value  = 44706
$ws.Cells.Item(800, 5).Value  = 13
$ws.Cells.Item(800, 6).Value  = 100112040
$ws.Cells.Item(800, 7).Value  = "Cilantro"
$ws.Cells.Item(800, 8).Value  = "Sin especificar"
$ws.Cells.Item(800, 9).Value  = "Primera"
$ws.Cells.Item(800, 10).Value = 350
$ws.Cells.Item(800, 11).Value = 7500
$ws.Cells.Item(800, 12).Value = 8000
$ws.Cells.Item(800, 13).Value = 7729
$ws.Cells.Item(800, 14).Value = "$/docena de atados"
$ws.Cells.Item(800, 15).Value = "Región Metropolitana"
$ws.Cells.Item(800, 16).Value = 2576
$ws.Cells.Item(800, 17).Value = 3
$ws.Cells.Item(800, 18).Value = "Hortaliza"
